# Refresh cryptos price/volume data (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.295.03'
$ws.Range("E2").Value = '  +0.02%  '

# Row 3
$ws.Range("D3").Value = '1.667.92'
$ws.Range("E3").Value = '  +0.21%  '

# Row 4
$ws.Range("D4").Value = '''1.008'
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$ws.Range("D5").Value = '''219.84'
$ws.Range("E5").Value = '  +0.49%  '

# Row 6
$ws.Range("D6").Value = '''0.5289'
$ws.Range("E6").Value = '  -0.51%  '

# Row 7
$ws.Range("D7").Value = '''1.009'
$ws.Range("E7").Value = '  -0.09%  '

# Row 8
$ws.Range("D8").Value = '''0.2648'
$ws.Range("E8").Value = '  +0.25%  '

# Row 9
$ws.Range("D9").Value = '''0.06365'
$ws.Range("E9").Value = '  -0.07%  '

# Row 10
$ws.Range("E10").Value = '  +1.89%  '

# Row 11
$ws.Range("D11").Value = '''0.07837'
$ws.Range("E11").Value = '  -0.19%  '

# Row 12
$ws.Range("D12").Value = '''4.523'
$ws.Range("E12").Value = '  -0.67%  '

# Row 13
$ws.Range("D13").Value = '1.660.68'
$ws.Range("E13").Value = '  -2.74%  '

# Row 14
$ws.Range("D14").Value = '1.895.68'
$ws.Range("E14").Value = '  +0.17%  '

# Row 15
$ws.Range("E15").Value = '  +1.24%  '

# Row 16
$ws.Range("D16").Value = '0.0₅8112'
$ws.Range("E16").Value = '  -1.14%  '

# Row 17
$ws.Range("D17").Value = '''65.70'
$ws.Range("E17").Value = '  -0.03%  '

# Row 18
$ws.Range("D18").Value = '26.303.94'
$ws.Range("E18").Value = '  +0.00%  '

# Row 19
$ws.Range("D19").Value = '''1.009'
$ws.Range("E19").Value = '  -0.01%  '

# Row 20
$ws.Range("D20").Value = '''4.719'
$ws.Range("E20").Value = '  +1.04%  '

# Row 21
$ws.Range("D21").Value = '''200.22'
$ws.Range("E21").Value = '  +3.91%  '

# Row 22
$ws.Range("E22").Value = '  +0.68%  '

# Row 23
$ws.Range("D23").Value = '''6.047'
$ws.Range("E23").Value = '  -0.25%  '

# Row 24
$ws.Range("D24").Value = '''1.010'
$ws.Range("E24").Value = '  -0.09%  '

# Row 25
$ws.Range("D25").Value = '''146.24'
$ws.Range("E25").Value = '  +0.71%  '

# Row 26
$ws.Range("D26").Value = '''0.1213'
$ws.Range("E26").Value = '  -0.72%  '

# Row 27
$ws.Range("D27").Value = '''7.237'
$ws.Range("E27").Value = '  -0.09%  '

# Row 28
$ws.Range("E28").Value = '  +0.09%  '

# Row 29
$ws.Range("D29").Value = '''1.530'
$ws.Range("E29").Value = '  +3.09%  '

# Row 30
$ws.Range("D30").Value = '''0.05908'
$ws.Range("E30").Value = '  +0.64%  '

# Row 31
$ws.Range("D31").Value = '''1.282'
$ws.Range("E31").Value = '  +0.19%  '

# Row 32
$ws.Range("D32").Value = '''3.516'
$ws.Range("E32").Value = '  -2.32%  '

# Row 33
$ws.Range("D33").Value = '''3.324'
$ws.Range("E33").Value = '  +0.53%  '

# Row 34
$ws.Range("E34").Value = '  -1.18%  '

# Row 35
$ws.Range("D35").Value = '''0.9632'

# Row 36
$ws.Range("D36").Value = '''2.822'
$ws.Range("E36").Value = '  +0.26%  '

# Row 37
$ws.Range("D37").Value = '''2.432'
$ws.Range("E37").Value = '  +0.20%  '

# Row 38
$ws.Range("D38").Value = '''0.5792'
$ws.Range("E38").Value = '  -0.33%  '

# Row 39
$ws.Range("D39").Value = '''0.01613'
$ws.Range("E39").Value = '  -0.02%  '

# Row 40
$ws.Range("D40").Value = '''5.953'
$ws.Range("E40").Value = '  +1.09%  '

# Row 41
$ws.Range("D41").Value = '1.075.62'
$ws.Range("E41").Value = '  +2.52%  '

# Row 42
$ws.Range("D42").Value = '''0.8575'
$ws.Range("E42").Value = '  +0.14%  '

# Row 43
$ws.Range("D43").Value = '''1.009'
$ws.Range("E43").Value = '  -0.02%  '

# Row 44
$ws.Range("D44").Value = '''102.90'
$ws.Range("E44").Value = '  -1.62%  '

# Row 45
$ws.Range("D45").Value = '1.806.32'
$ws.Range("E45").Value = '  +0.04%  '

# Row 46
$ws.Range("D46").Value = '''58.46'
$ws.Range("E46").Value = '  +1.93%  '

# Row 47
$ws.Range("E47").Value = '  -0.21%  '

# Row 48
$ws.Range("D48").Value = '''0.4414'
$ws.Range("E48").Value = '  +0.88%  '

# Row 49
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₈104'
$ws.Range("E49").Value = '  -2.54%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''8.066'
$ws.Range("E50").Value = '  +1.37%  '

# Row 51
$ws.Range("D51").Value = '''0.05144'
$ws.Range("E51").Value = '  -0.38%  '

